$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C) column date for all data rows (2-515) to 2023-09-19 (serial 45188)
$ws.Range("C2:C515").Value = "2023-09-19"

# Row 2 specific updates: Fridlysta (H2) 10 -> 11, Alla arter (Q2) 32 -> 33
$ws.Range("H2").Value = 11
$ws.Range("Q2").Value = 33

# Insert "Skogsrör" into the species list (Artnamn, R2) between "Tvåblad" and "Brudsporre"
$ws.Range("R2").Value = "Fjällfotad musseron`r`nGrantickeporing`r`nKnärot`r`nLäderdoftande fingersvamp`r`nSkäggvaxskivling`r`nSmalskaftslav`r`nSprickporing`r`nBarrviolspindling`r`nFlattoppad klubbsvamp`r`nGarnlav`r`nGranticka`r`nGultoppig fingersvamp`r`nLuddfingersvamp`r`nLunglav`r`nRödbrun klubbdyna`r`nSkogsfru`r`nSpillkråka`r`nTretåig hackspett`r`nUllticka`r`nVitterspindling`r`nÄggvaxskivling`r`nDropptaggsvamp`r`nFinbräken`r`nGuckusko`r`nSkinnlav`r`nSkogsknipprot`r`nSvavelriska`r`nTrådticka`r`nTvåblad`r`nSkogsrör`r`nBrudsporre`r`nFläcknycklar`r`nÄngsnycklar"

# Keep the row height consistent with the rest of the sheet (avoid auto row-height growth
# from the extra wrapped line that was just added to the cell above).
$ws.Rows.Item(2).RowHeight = 15
